# Auto-generated: update FFXIV market price / profit data cells across all 8 sheets
# per the scheduled runner's refreshed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 58.333332
$ws.Range("I12").Value = 58.333332
$ws.Range("K12").Value = 58.333332
$ws.Range("M12").Value = 111.666668
# Row 34
$ws.Range("H34").Value = 3001
$ws.Range("I34").Value = 3001
$ws.Range("K34").Value = 3001
$ws.Range("M34").Value = -2798
# Row 36
$ws.Range("H36").Value = 3001
$ws.Range("I36").Value = 3001
$ws.Range("K36").Value = 3001
$ws.Range("M36").Value = -2286
# Row 70
$ws.Range("H70").Value = 1619.7
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 1466.3334
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 4399.0002
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -4939.0002
# Row 73
$ws.Range("H73").Value = 1619.7
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 1466.3334
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 4399.0002
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -6271.0002
# Row 132
$ws.Range("H132").Value = 2075.1538
$ws.Range("I132").Value = 1713.8889
$ws.Range("J132").Value = 2888
$ws.Range("K132").Value = 5141.6667
$ws.Range("L132").Value = 8664
$ws.Range("M132").Value = -2611.6667
$ws.Range("N132").Value = -13724
# Row 133
$ws.Range("H133").Value = 77377.89999999999
$ws.Range("J133").Value = 77377.89999999999
$ws.Range("L133").Value = 77377.89999999999
$ws.Range("N133").Value = -87497.89999999999
# Row 134
$ws.Range("H134").Value = 98988.57000000001
$ws.Range("J134").Value = 98988.57000000001
$ws.Range("L134").Value = 98988.57000000001
$ws.Range("N134").Value = -109128.57
# Row 136
$ws.Range("H136").Value = 77977.14
$ws.Range("J136").Value = 77977.14
$ws.Range("L136").Value = 77977.14
$ws.Range("N136").Value = -88177.14
# Row 139
$ws.Range("H139").Value = 74505.22
$ws.Range("J139").Value = 74505.22
$ws.Range("L139").Value = 74505.22
$ws.Range("N139").Value = -84785.22
# Row 140
$ws.Range("H140").Value = 49649.727
$ws.Range("J140").Value = 49437.555
$ws.Range("L140").Value = 49437.555
$ws.Range("N140").Value = -59797.555

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
# Row 49
$ws.Range("H49").Value = 39900
$ws.Range("J49").Value = 39900
$ws.Range("L49").Value = 39900
$ws.Range("N49").Value = -40420
# Row 74
$ws.Range("H74").Value = 35437.832
$ws.Range("I74").Value = 44793.176
$ws.Range("K74").Value = 44793.176
$ws.Range("M74").Value = -43919.176
# Row 77
$ws.Range("H77").Value = 35437.832
$ws.Range("I77").Value = 44793.176
$ws.Range("K77").Value = 223965.88
$ws.Range("M77").Value = -219597.88
# Row 122
$ws.Range("H122").Value = 2483.1667
$ws.Range("I122").Value = 2483.1667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7449.500100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4999.500100000001
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 2650.8845
$ws.Range("I132").Value = 2377.7896
$ws.Range("K132").Value = 7133.3688
$ws.Range("M132").Value = -4603.3688

$ws = $wb.Worksheets.Item("BSM")
# Row 53
$ws.Range("H53").Value = 38964
$ws.Range("J53").Value = 38964
$ws.Range("L53").Value = 38964
$ws.Range("N53").Value = -40112
# Row 107
$ws.Range("H107").Value = 14288279
$ws.Range("I107").Value = 25002256
$ws.Range("J107").Value = 2977
$ws.Range("K107").Value = 25002256
$ws.Range("L107").Value = 2977
$ws.Range("M107").Value = -25000336
$ws.Range("N107").Value = -6817
# Row 109
$ws.Range("H109").Value = 78282.71000000001
$ws.Range("J109").Value = 78282.71000000001
$ws.Range("L109").Value = 78282.71000000001
$ws.Range("N109").Value = -81056.71000000001
# Row 132
$ws.Range("H132").Value = 26321.695
$ws.Range("J132").Value = 26321.695
$ws.Range("L132").Value = 26321.695
$ws.Range("N132").Value = -36441.695
# Row 134
$ws.Range("H134").Value = 4419.9355
$ws.Range("I134").Value = 2724.0952
$ws.Range("K134").Value = 8172.285600000001
$ws.Range("M134").Value = -5637.285600000001
# Row 135
$ws.Range("H135").Value = 99995.8
$ws.Range("J135").Value = 99995.8
$ws.Range("L135").Value = 99995.8
$ws.Range("N135").Value = -110135.8
# Row 138
$ws.Range("H138").Value = 77893.5
$ws.Range("J138").Value = 77893.5
$ws.Range("L138").Value = 77893.5
$ws.Range("N138").Value = -88173.5
# Row 140
$ws.Range("H140").Value = 93496
$ws.Range("J140").Value = 93496
$ws.Range("L140").Value = 93496
$ws.Range("N140").Value = -103856

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 466.66666
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -1700
# Row 58
$ws.Range("H58").Value = 1437.921
$ws.Range("I58").Value = 1290.1936
$ws.Range("K58").Value = 1290.1936
$ws.Range("M58").Value = -1087.1936
# Row 62
$ws.Range("H62").Value = 2899.5
$ws.Range("I62").Value = 2999.4443
$ws.Range("K62").Value = 2999.4443
$ws.Range("M62").Value = -2375.4443
# Row 65
$ws.Range("H65").Value = 2899.5
$ws.Range("I65").Value = 2999.4443
$ws.Range("K65").Value = 14997.2215
$ws.Range("M65").Value = -11877.2215
# Row 134
$ws.Range("H134").Value = 35909.9
$ws.Range("I134").Value = 2584.5217
$ws.Range("K134").Value = 7753.5651
$ws.Range("M134").Value = -5218.5651
# Row 136
$ws.Range("H136").Value = 1437.921
$ws.Range("I136").Value = 1290.1936
$ws.Range("K136").Value = 3870.5808
$ws.Range("M136").Value = -1320.5808
# Row 138
$ws.Range("H138").Value = 54353.332
$ws.Range("J138").Value = 54897.5
$ws.Range("L138").Value = 54897.5
$ws.Range("N138").Value = -65177.5

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 476.66666
$ws.Range("J75").Value = 490
$ws.Range("L75").Value = 1470
$ws.Range("N75").Value = -3466
# Row 78
$ws.Range("H78").Value = 476.66666
$ws.Range("J78").Value = 490
$ws.Range("L78").Value = 4410
$ws.Range("N78").Value = -14394
# Row 103
$ws.Range("H103").Value = 952.25
$ws.Range("I103").Value = 936.3333
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 2808.9999
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = -1929.9999
$ws.Range("N103").Value = -4758
# Row 129
$ws.Range("I129").Value = 371.85715
$ws.Range("J129").Value = 333333340
$ws.Range("K129").Value = 1115.57145
$ws.Range("L129").Value = 1000000020
$ws.Range("M129").Value = 3884.42855
$ws.Range("N129").Value = -1000010020

$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 13083828
$ws.Range("J33").Value = 13340814
$ws.Range("L33").Value = 13340814
$ws.Range("N33").Value = -13341318
# Row 49
$ws.Range("H49").Value = 20711.8
$ws.Range("J49").Value = 20711.8
$ws.Range("L49").Value = 20711.8
$ws.Range("N49").Value = -21079.8
# Row 93
$ws.Range("H93").Value = 16404.9
$ws.Range("J93").Value = 16404.9
$ws.Range("L93").Value = 16404.9
$ws.Range("N93").Value = -20148.9
# Row 109
$ws.Range("H109").Value = 17499.7
$ws.Range("J109").Value = 18555.223
$ws.Range("L109").Value = 18555.223
$ws.Range("N109").Value = -20635.223
# Row 140
$ws.Range("H140").Value = 94552.37
$ws.Range("J140").Value = 94957.60000000001
$ws.Range("L140").Value = 94957.60000000001
$ws.Range("N140").Value = -105317.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3029.2307
$ws.Range("I7").Value = 2098.75
$ws.Range("K7").Value = 2098.75
$ws.Range("M7").Value = -1986.75
# Row 40
$ws.Range("H40").Value = 7940151.5
$ws.Range("I40").Value = 3125.75
$ws.Range("K40").Value = 3125.75
$ws.Range("M40").Value = -2989.75
# Row 68
$ws.Range("H68").Value = 2780.1
$ws.Range("I68").Value = 2850.125
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 2850.125
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -2101.125
$ws.Range("N68").Value = -3998
# Row 71
$ws.Range("H71").Value = 2780.1
$ws.Range("I71").Value = 2850.125
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 14250.625
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -10506.625
$ws.Range("N71").Value = -19988
# Row 126
$ws.Range("H126").Value = 3029.2307
$ws.Range("I126").Value = 2098.75
$ws.Range("K126").Value = 6296.25
$ws.Range("M126").Value = -3826.25

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 116
$ws.Range("H116").Value = 77387
$ws.Range("J116").Value = 77387
$ws.Range("L116").Value = 77387
$ws.Range("N116").Value = -86565
# Row 133
$ws.Range("H133").Value = 50681
$ws.Range("J133").Value = 50681
$ws.Range("L133").Value = 50681
$ws.Range("N133").Value = -60801
